# Highlight the three "buddy settings" requirement bullets (admin related
# permissions) in yellow, matching the author's note about the admin
# add/delete settings.  We match on the start of each bullet's text and
# skip any occurrence that lives inside a table (the summary table further
# down the document repeats some of this text but must stay untouched).

$d = $word.ActiveDocument

$targets = @(
    "אוכל לנהל הרשאות של הדרייב",
    "אוכל לנהל הרשאות למשתמשי הדרייב",
    "אוכל למנות תתי אדמין ולתת להם הרשאות"
)

$wdWithInTable = 12

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $rng = $para.Range

    if ($rng.Information($wdWithInTable)) {
        continue
    }

    $text = $rng.Text
    foreach ($target in $targets) {
        if ($text.StartsWith($target)) {
            # Setting Font.HighlightColorIndex (rather than the Range's own
            # HighlightColorIndex) applies the highlight both to every run
            # in the paragraph AND to the paragraph mark's run properties,
            # mirroring how the rest of the document was highlighted.
            $rng.Font.HighlightColorIndex = 7   # wdYellow
            break
        }
    }
}
